# Auto-generated edit script applying the Cerberus_Profits.xlsx commit diff.
# Each write reproduces one <c> value change (or add/remove) from the diff,
# grouped by worksheet (tab order: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3313.875  # was 3478.875
$ws.Range("J40").Value = 3369.75  # was 3699.75
$ws.Range("L40").Value = 3369.75  # was 3699.75
$ws.Range("N40").Value = -3719.75  # was -4049.75
$ws.Range("H51").Value = 12107.223  # was 36745.75
$ws.Range("I51").Value = 14300.2  # was 52800.4
$ws.Range("J51").Value = 9366  # was 9988
$ws.Range("K51").Value = 14300.2  # was 52800.4
$ws.Range("L51").Value = 9366  # was 9988
$ws.Range("M51").Value = -13816.2  # was -52316.4
$ws.Range("N51").Value = -10334  # was -10956
$ws.Range("H74").Value = 7315.472  # was 7388.8
$ws.Range("I74").Value = 6747.933  # was 6890.7144
$ws.Range("K74").Value = 6747.933  # was 6890.7144
$ws.Range("M74").Value = -5811.933  # was -5954.7144
$ws.Range("H77").Value = 7315.472  # was 7388.8
$ws.Range("I77").Value = 6747.933  # was 6890.7144
$ws.Range("K77").Value = 33739.665  # was 34453.572
$ws.Range("M77").Value = -29059.665  # was -29773.572
$ws.Range("H116").Value = 7111.5884  # was 7111.647
$ws.Range("I116").Value = 8658.916999999999  # was 8659
$ws.Range("K116").Value = 8658.916999999999  # was 8659
$ws.Range("M116").Value = -5216.916999999999  # was -5217
$ws.Range("H135").Value = 27122.357  # was 19506.7
$ws.Range("I135").Value = 2029.375  # was 1917
$ws.Range("J135").Value = 60579.668  # was 52173.285
$ws.Range("K135").Value = 18264.375  # was 17253
$ws.Range("L135").Value = 545217.012  # was 469559.5650000001
$ws.Range("M135").Value = -15729.375  # was -14718
$ws.Range("N135").Value = -550287.012  # was -474629.5650000001
$ws.Range("H137").Value = 3602.6428  # was 3596
$ws.Range("I137").Value = 1745.25  # was 1746.125
$ws.Range("J137").Value = 6079.1665  # was 6062.5
$ws.Range("K137").Value = 5235.75  # was 5238.375
$ws.Range("L137").Value = 18237.4995  # was 18187.5
$ws.Range("M137").Value = -2685.75  # was -2688.375
$ws.Range("N137").Value = -23337.4995  # was -23287.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 48999.5  # was 7500
$ws.Range("J106").Value = 48999.5  # was 7500
$ws.Range("L106").Value = 48999.5  # was 7500
$ws.Range("N106").Value = -51523.5  # was -10024
$ws.Range("H122").Value = 3391.2666  # was 11331.4
$ws.Range("I122").Value = 3765.5557  # was 16999.111
$ws.Range("K122").Value = 11296.6671  # was 50997.333
$ws.Range("M122").Value = -8846.667099999999  # was -48547.333
$ws.Range("H132").Value = 3215.7585  # was 3294.9285
$ws.Range("I132").Value = 3133.348  # was 3230.3635
$ws.Range("K132").Value = 9400.044  # was 9691.0905
$ws.Range("M132").Value = -6870.044  # was -7161.0905
$ws.Range("H133").Value = 316666.34  # was 250000
$ws.Range("J133").Value = 449999.5  # was 450000
$ws.Range("L133").Value = 449999.5  # was 450000
$ws.Range("N133").Value = -455059.5  # was -455060

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 52307.11  # was 46775.89
$ws.Range("J58").Value = 52307.11  # was 46775.89
$ws.Range("L58").Value = 52307.11  # was 46775.89
$ws.Range("N58").Value = -52895.11  # was -47363.89
$ws.Range("H60").Value = 16141.571  # was 16665.5
$ws.Range("J60").Value = 16141.571  # was 16665.5
$ws.Range("L60").Value = 16141.571  # was 16665.5
$ws.Range("N60").Value = -17339.571  # was -17863.5
$ws.Range("H99").Value = 2111.5334  # was 2109.6875
$ws.Range("I99").Value = 1995.3846  # was 2001.5714
$ws.Range("K99").Value = 1995.3846  # was 2001.5714
$ws.Range("M99").Value = -497.3846000000001  # was -503.5714
$ws.Range("H105").Value = 3191.1  # was 3257.6667
$ws.Range("I105").Value = 2822.9167  # was 2956.3914
$ws.Range("J105").Value = 4663.8335  # was 4247.5713
$ws.Range("K105").Value = 2822.9167  # was 2956.3914
$ws.Range("L105").Value = 4663.8335  # was 4247.5713
$ws.Range("M105").Value = -1075.9167  # was -1209.3914
$ws.Range("N105").Value = -8157.8335  # was -7741.5713

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2695.5833  # was 2638.348
$ws.Range("J16").Value = 4794.4  # was 4881.3335
$ws.Range("L16").Value = 4794.4  # was 4881.3335
$ws.Range("N16").Value = -5368.4  # was -5455.3335
$ws.Range("H22").Value = 461.7143  # was 401.33334
$ws.Range("I22").Value = 434  # was 395
$ws.Range("J22").Value = 531  # was 414
$ws.Range("K22").Value = 434  # was 395
$ws.Range("L22").Value = 531  # was 414
$ws.Range("M22").Value = -84  # was -45
$ws.Range("N22").Value = -1231  # was -1114
$ws.Range("H113").Value = 2695.5833  # was 2638.348
$ws.Range("J113").Value = 4794.4  # was 4881.3335
$ws.Range("L113").Value = 4794.4  # was 4881.3335
$ws.Range("N113").Value = -9134.4  # was -9221.333500000001
$ws.Range("H132").Value = 2262.0454  # was 2330.476
$ws.Range("I132").Value = 2564.9412  # was 2577.7646
$ws.Range("J132").Value = 1232.2  # was 1279.5
$ws.Range("K132").Value = 7694.823600000001  # was 7733.293799999999
$ws.Range("L132").Value = 3696.6  # was 3838.5
$ws.Range("M132").Value = -5164.823600000001  # was -5203.293799999999
$ws.Range("N132").Value = -8756.6  # was -8898.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1449.8889  # was 1567.7778
$ws.Range("I7").Value = 189.75  # was 243.33333
$ws.Range("J7").Value = 2458  # was 2230
$ws.Range("K7").Value = 569.25  # was 729.99999
$ws.Range("L7").Value = 7374  # was 6690
$ws.Range("M7").Value = -457.25  # was -617.99999
$ws.Range("N7").Value = -7598  # was -6914
$ws.Range("H11").Value = 3135698.8  # was 2743736.5
$ws.Range("I11").Value = 3658265.2  # was 3135656.2
$ws.Range("K11").Value = 10974795.6  # was 9406968.600000001
$ws.Range("M11").Value = -10974655.6  # was -9406828.600000001
$ws.Range("H14").Value = 943.6667  # was 842.5714
$ws.Range("I14").Value = 943.6667  # was 842.5714
$ws.Range("K14").Value = 2831.0001  # was 2527.7142
$ws.Range("M14").Value = -2658.0001  # was -2354.7142
$ws.Range("H40").Value = 4625.75  # was 5279.7144
$ws.Range("J40").Value = 5513.25  # was 7335
$ws.Range("L40").Value = 22053  # was 29340
$ws.Range("N40").Value = -22191  # was -29478
$ws.Range("H80").Value = 10000  # was 0
$ws.Range("J80").Value = 10000  # was 0
$ws.Range("L80").Value = 30000  # was 0
$ws.Range("N80").Value = -31872  # new cell
$ws.Range("H83").Value = 10000  # was 0
$ws.Range("J83").Value = 10000  # was 0
$ws.Range("L83").Value = 90000  # was 0
$ws.Range("N83").Value = -99360  # new cell
$ws.Range("H92").Value = 3465.3333  # was 3509.7778
$ws.Range("J92").Value = 3110.75  # was 3160.75
$ws.Range("L92").Value = 9332.25  # was 9482.25
$ws.Range("N92").Value = -11828.25  # was -11978.25
$ws.Range("H114").Value = 8510.333000000001  # was 7467.75
$ws.Range("I114").Value = 0  # was 840
$ws.Range("J114").Value = 8510.333000000001  # was 9677
$ws.Range("K114").Value = 0  # was 2520
$ws.Range("L114").Value = 25530.999  # was 29031
$ws.Range("M114").ClearContents()  # was 734, removed
$ws.Range("N114").Value = -32038.999  # was -35539
$ws.Range("H117").Value = 3027.65  # was 3068.762
$ws.Range("I117").Value = 1571.3334  # was 1805.6
$ws.Range("J117").Value = 3651.7856  # was 3463.5
$ws.Range("K117").Value = 4714.0002  # was 5416.799999999999
$ws.Range("L117").Value = 10955.3568  # was 10390.5
$ws.Range("M117").Value = -1272.0002  # was -1974.799999999999
$ws.Range("N117").Value = -17839.3568  # was -17274.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3748.9614  # was 3771.5386
$ws.Range("I132").Value = 3748.9614  # was 3771.5386
$ws.Range("K132").Value = 11246.8842  # was 11314.6158
$ws.Range("M132").Value = -8716.8842  # was -8784.6158

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3149.7778  # was 3063.2222
$ws.Range("I22").Value = 821  # was 735.5
$ws.Range("J22").Value = 3440.875  # was 3728.2856
$ws.Range("K22").Value = 821  # was 735.5
$ws.Range("L22").Value = 3440.875  # was 3728.2856
$ws.Range("M22").Value = -526  # was -440.5
$ws.Range("N22").Value = -4030.875  # was -4318.2856
$ws.Range("H27").Value = 3149.7778  # was 3063.2222
$ws.Range("I27").Value = 821  # was 735.5
$ws.Range("J27").Value = 3440.875  # was 3728.2856
$ws.Range("K27").Value = 821  # was 735.5
$ws.Range("L27").Value = 3440.875  # was 3728.2856
$ws.Range("M27").Value = -714  # was -628.5
$ws.Range("N27").Value = -3654.875  # was -3942.2856
$ws.Range("H100").Value = 2498.5  # was 2498.5833
$ws.Range("I100").Value = 1604.6666  # was 1604.8334
$ws.Range("K100").Value = 1604.6666  # was 1604.8334
$ws.Range("M100").Value = -1063.6666  # was -1063.8334

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 120000  # was 94999.5
$ws.Range("J46").Value = 120000  # was 94999.5
$ws.Range("L46").Value = 120000  # was 94999.5
$ws.Range("N46").Value = -120462  # was -95461.5
$ws.Range("H62").Value = 28125.5  # was 36667.332
$ws.Range("I62").Value = 7500.6665  # was 10001
$ws.Range("K62").Value = 7500.6665  # was 10001
$ws.Range("M62").Value = -6876.6665  # was -9377
$ws.Range("H65").Value = 28125.5  # was 36667.332
$ws.Range("I65").Value = 7500.6665  # was 10001
$ws.Range("K65").Value = 37503.3325  # was 50005
$ws.Range("M65").Value = -34383.3325  # was -46885
$ws.Range("I81").Value = 4437  # was 4718.2856
$ws.Range("J81").Value = 53001  # was 35500.332
$ws.Range("K81").Value = 8874  # was 9436.5712
$ws.Range("L81").Value = 106002  # was 71000.664
$ws.Range("M81").Value = -7813  # was -8375.5712
$ws.Range("N81").Value = -108124  # was -73122.664
$ws.Range("I84").Value = 4437  # was 4718.2856
$ws.Range("J84").Value = 53001  # was 35500.332
$ws.Range("K84").Value = 44370  # was 47182.856
$ws.Range("L84").Value = 530010  # was 355003.32
$ws.Range("M84").Value = -39066  # was -41878.856
$ws.Range("N84").Value = -540618  # was -365611.32
$ws.Range("H132").Value = 13516774  # was 13892229
$ws.Range("I132").Value = 18871064  # was 19611090
$ws.Range("K132").Value = 56613192  # was 58833270
$ws.Range("M132").Value = -56610662  # was -58830740
$ws.Range("H134").Value = 120000  # was 94999.5
$ws.Range("J134").Value = 120000  # was 94999.5
$ws.Range("L134").Value = 360000  # was 284998.5
$ws.Range("N134").Value = -365070  # was -290068.5

